$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at E:F (old E/F shift right to G/H, keeping their
# formulas/values/styles intact; new E/F inherit the style of column D).
$ws.Columns("E:F").Insert()

# New header row (row 2) values for the inserted columns.
$ws.Range("E2").Value = "price"
$ws.Range("F2").Value = "max"

# New numeric data for rows 3-12.
$ws.Range("E3").Value = 123
$ws.Range("F3").Value = 100

$ws.Range("E4").Value = 321
$ws.Range("F4").Value = 50

$ws.Range("E5").Value = 456
$ws.Range("F5").Value = 99

$ws.Range("E6").Value = 4895
$ws.Range("F6").Value = 200

$ws.Range("E7").Value = 21546
$ws.Range("F7").Value = 200

$ws.Range("E8").Value = 45.54
$ws.Range("F8").Value = 99

$ws.Range("E9").Value = 453.54
$ws.Range("F9").Value = 50

$ws.Range("E10").Value = 5.9
$ws.Range("F10").Value = 50

$ws.Range("E11").Value = 9.99
$ws.Range("F11").Value = 100

$ws.Range("E12").Value = 1003.5
$ws.Range("F12").Value = 1

# Cosmetic changes: default column width, default row height, selection.
$ws.Columns("A:H").ColumnWidth = 17.4

$ws.Application.GoTo($ws.Range("E13"))

Write-Host "done"
